# "Done due date set up"
# Rows 18 and 19 (the GPB/GAD AR due date feature and the monitoring
# table) are now finished: their Status column flips from PENDING to OK
# and a Remark of "Done" is filled in, matching the rest of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "OK"
$ws.Range("E18").Value = "Done"

$ws.Range("D19").Value = "OK"
$ws.Range("E19").Value = "Done"

# Keep the active selection in sync with the newly completed row, the
# same way Excel nudges it forward as you work down the sheet.
$ws.Range("D19").Select()
